$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 173, shifting existing rows 173:285 down to 174:286
$ws.Rows(173).Insert()

# Populate the new row 173 with the new record
$ws.Cells.Item(173, 1).Value = 8
$ws.Cells.Item(173, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(173, 3).Value = "Coquimbo"
$ws.Cells.Item(173, 4).Value = (Get-Date -Year 2022 -Month 2 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(173, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(173, 5).Value = 4
$ws.Cells.Item(173, 6).Value = 100114013
$ws.Cells.Item(173, 7).Value = "Zanahoria"
$ws.Cells.Item(173, 8).Value = "Sin especificar"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 600
$ws.Cells.Item(173, 11).Value = 5500
$ws.Cells.Item(173, 12).Value = 6000
$ws.Cells.Item(173, 13).Value = 5750
$ws.Cells.Item(173, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(173, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(173, 16).Value = 288
$ws.Cells.Item(173, 17).Value = 20
$ws.Cells.Item(173, 18).Value = "Hortaliza"
